$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2023-10-08 10:42:34", "hatespeech", "setting1", 2, 3844.990900800918),
    @("2023-10-08 10:42:35", "hatespeech", "def",      2, 2893.15572850852),
    @("2023-10-08 10:42:35", "hatespeech", "setting2", 2, 2354.425326036787),
    @("2023-10-08 10:42:35", "hatespeech", "setting3", 2, 3783.355944707603),
    @("2023-10-08 10:42:35", "hatespeech", "setting5", 2, 2298.094793551396),
    @("2023-10-08 10:42:35", "hatespeech", "setting4", 2, 2868.605898795541)
)

$startRow = 27
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
